$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (the "Förändrad" / date-changed column) for rows 2 through 39
# from serial date 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C39").Value = 45202
